$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "EQD-6,EQN-11" -> "EQD-7,EQN-11" (row 9, B and C columns)
$ws.Range("B9").Value = "EQD-7,EQN-11"
$ws.Range("C9").Value = "EQD-7,EQN-11"

# Add new requisito row 26, mirroring row 25's formatting exactly (wrap text,
# top-aligned; B = default/black font, C = red font) by copying the format
# from row 25's cells instead of toggling properties one-by-one (which would
# otherwise leave stray unused style records behind).
$ws.Range("B26").Value = "LOQ4082 -  Corrosão  (Requisito fraco)`n"
$ws.Range("B25").Copy()
$ws.Range("B26").PasteSpecial(-4122)

$ws.Range("C26").Value = "LOQ4082 -  Corrosão  (Requisito fraco)`n"
$ws.Range("C25").Copy()
$ws.Range("C26").PasteSpecial(-4122)

$ws.Rows.Item(26).RowHeight = 30
